# Jackie - xlxx update
#
# Adds Approach / Pass Criteria / Fail Criteria text for the "Timeline",
# "Symptom Checker" and "Alert!" rows (rows 3, 4, 5 respectively) to the
# Features_Security_Testing sheet, and grows those rows to fit the new
# wrapped text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Timeline ---------------------------------------------------
$ws.Range("G3").Value = "We will test to ensure the images and content that are placed within the history timeline display properly"
$ws.Range("H3").Value = "The timeline slider display in the proper order withotu any error"
$ws.Range("I3").Value = "The images and content do not advnace through the series as they should, display them out of the proper sequence, without the disiignated interval, any of the image are not properly called from the database."

# --- Row 5: Alert! -------------------------------------------------------
$ws.Range("G5").Value = "We will ensure the message that are place within the alert display properly"
$ws.Range("H5").Value = "The alert display in the proper order without any error"
$ws.Range("I5").Value = "The message do not advnace through the series as it should, display it out of the proper sequence, without the disiignated interval, any of the message are not properly called from the database."

# --- Row 4: Symptom Checker ----------------------------------------------
$ws.Range("G4").Value = "We will input a few symptom in to the database and the appropriate symptom will display base on the user's f&q input"
$ws.Range("H4").Value = "The appropriate symptom display successfully base on the user's input"
$ws.Range("I4").Value = "The symptom fail to display"

# Match the existing wrap-text style used by the rest of column D/G/H/I.
$ws.Range("G3:I3").WrapText = $true
$ws.Range("G4:I4").WrapText = $true
$ws.Range("G5:I5").WrapText = $true

# Grow the rows so the new wrapped text is fully visible.
$ws.Rows(3).RowHeight = 120
$ws.Rows(4).RowHeight = 45
$ws.Rows(5).RowHeight = 120

# Move the active selection, matching the saved view state.
$ws.Range("K4").Select() | Out-Null
